# Update workbook "26064948_Rio Negro.xlsx" with the latest daily readings
# (data through 2020-08-23 / serial 44066), per commit
# "Atualização dos gráficos 26082020".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Mensal": last monthly summary row (row 14) now reflects the
# newest date / values instead of the previous month-end snapshot.
# ---------------------------------------------------------------------
$wsMensal = $wb.Worksheets.Item("Mensal")

$wsMensal.Range("A14").Value = 44066
$wsMensal.Range("B14").Value = 124.21
$wsMensal.Range("D14").Value = 10.17
# C14 (112.75) is unchanged.

# ---------------------------------------------------------------------
# Sheet "Diario": append seven new daily rows (384-390), one for each
# day from 2020-08-17 through 2020-08-23.
# ---------------------------------------------------------------------
$wsDiario = $wb.Worksheets.Item("Diario")

$newRows = @(
    @{ Row = 384; Data = 44060; Obs = 259.41; Ref = 112.75; Anom = 130.08 },
    @{ Row = 385; Data = 44061; Obs = 226.34; Ref = 112.75; Anom = 100.75 },
    @{ Row = 386; Data = 44062; Obs = 218.85; Ref = 112.75; Anom = 94.1 },
    @{ Row = 387; Data = 44063; Obs = 231.61; Ref = 112.75; Anom = 105.42 },
    @{ Row = 388; Data = 44064; Obs = 246.3;  Ref = 112.75; Anom = 118.45 },
    @{ Row = 389; Data = 44065; Obs = 253.17; Ref = 112.75; Anom = 124.54 },
    @{ Row = 390; Data = 44066; Obs = 237.53; Ref = 112.75; Anom = 110.67 }
)

foreach ($r in $newRows) {
    $rowNum = $r.Row
    $wsDiario.Range("A$rowNum").Value = $r.Data
    $wsDiario.Range("B$rowNum").Value = $r.Obs
    $wsDiario.Range("C$rowNum").Value = $r.Ref
    $wsDiario.Range("D$rowNum").Value = $r.Anom

    # Copy the date-column formatting (number format, bold font, border)
    # from the previous row so the new date cell matches the existing
    # style (s="2") instead of getting the default style.
    $prevRow = $rowNum - 1
    $wsDiario.Range("A$prevRow").Copy()
    $wsDiario.Range("A$rowNum").PasteSpecial(-4122) # xlPasteFormats
}

$excel.CutCopyMode = 0
